$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 11.24958966666667
$ws.Range("H2").Value = 33.748769
$ws.Range("I2").Value = 0.04815412300202451
$ws.Range("J2").Value = 0.04815412300202451
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 10.779612
$ws.Range("N2").Value = 32.338836
$ws.Range("O2").Value = 0.1321092878737708
$ws.Range("P2").Value = 0.1321092878737708
$ws.Range("Q2").Value = 121.266211765876
$ws.Range("R2").Value = 1091.395905892884
$ws.Range("S2").Value = 0.006361606897983426
$ws.Range("T2").Value = 0.006361606897983423
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 11.24958966666667
$ws.Range("H3").Value = 33.748769
$ws.Range("I3").Value = 0.04815412300202451
$ws.Range("J3").Value = 0.04815412300202451
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 38.54369466666667
$ws.Range("N3").Value = 115.631084
$ws.Range("O3").Value = 0.4723713668393066
$ws.Range("P3").Value = 0.4723713668393065
$ws.Range("Q3").Value = 433.6007492372884
$ws.Range("R3").Value = 3902.406743135596
$ws.Range("S3").Value = 0.02274662890141441
$ws.Range("T3").Value = 0.02274662890141441
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 11.24958966666667
$ws.Range("H4").Value = 33.748769
$ws.Range("I4").Value = 0.04815412300202451
$ws.Range("J4").Value = 0.04815412300202451
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 12.62567333333333
$ws.Range("N4").Value = 37.87702
$ws.Range("O4").Value = 0.1547336502458089
$ws.Range("P4").Value = 0.1547336502458089
$ws.Range("Q4").Value = 142.0336442653755
$ws.Range("R4").Value = 1278.30279838838
$ws.Range("S4").Value = 0.007451063226488924
$ws.Range("T4").Value = 0.007451063226488921
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 11.24958966666667
$ws.Range("H5").Value = 33.748769
$ws.Range("I5").Value = 0.04815412300202451
$ws.Range("J5").Value = 0.04815412300202451
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 19.64719066666667
$ws.Range("N5").Value = 58.94157200000001
$ws.Range("O5").Value = 0.2407856950411137
$ws.Range("P5").Value = 0.2407856950411137
$ws.Range("Q5").Value = 221.0228331027631
$ws.Range("R5").Value = 1989.205497924868
$ws.Range("S5").Value = 0.01159482397613775
$ws.Range("T5").Value = 0.01159482397613775
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 208.2711843333333
$ws.Range("H6").Value = 624.813553
$ws.Range("I6").Value = 0.8915095150431698
$ws.Range("J6").Value = 0.8915095150431697
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 10.779612
$ws.Range("N6").Value = 32.338836
$ws.Range("O6").Value = 0.1321092878737708
$ws.Range("P6").Value = 0.1321092878737708
$ws.Range("Q6").Value = 2245.082557893812
$ws.Range("R6").Value = 20205.74302104431
$ws.Range("S6").Value = 0.1177766871650439
$ws.Range("T6").Value = 0.1177766871650439
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 208.2711843333333
$ws.Range("H7").Value = 624.813553
$ws.Range("I7").Value = 0.8915095150431698
$ws.Range("J7").Value = 0.8915095150431697
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 38.54369466666667
$ws.Range("N7").Value = 115.631084
$ws.Range("O7").Value = 0.4723713668393066
$ws.Range("P7").Value = 0.4723713668393065
$ws.Range("Q7").Value = 8027.54093680905
$ws.Range("R7").Value = 72247.86843128144
$ws.Range("S7").Value = 0.4211235681711895
$ws.Range("T7").Value = 0.4211235681711893
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 208.2711843333333
$ws.Range("H8").Value = 624.813553
$ws.Range("I8").Value = 0.8915095150431698
$ws.Range("J8").Value = 0.8915095150431697
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 12.62567333333333
$ws.Range("N8").Value = 37.87702
$ws.Range("O8").Value = 0.1547336502458089
$ws.Range("P8").Value = 0.1547336502458089
$ws.Range("Q8").Value = 2629.563938139117
$ws.Range("R8").Value = 23666.07544325206
$ws.Range("S8").Value = 0.1379465214915006
$ws.Range("T8").Value = 0.1379465214915005
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 208.2711843333333
$ws.Range("H9").Value = 624.813553
$ws.Range("I9").Value = 0.8915095150431698
$ws.Range("J9").Value = 0.8915095150431697
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 19.64719066666667
$ws.Range("N9").Value = 58.94157200000001
$ws.Range("O9").Value = 0.2407856950411137
$ws.Range("P9").Value = 0.2407856950411137
$ws.Range("Q9").Value = 4091.94366896948
$ws.Range("R9").Value = 36827.49302072532
$ws.Range("S9").Value = 0.2146627382154359
$ws.Range("T9").Value = 0.2146627382154358
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.262842333333333
$ws.Range("H10").Value = 9.788527
$ws.Range("I10").Value = 0.01396667040408609
$ws.Range("J10").Value = 0.01396667040408609
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 10.779612
$ws.Range("N10").Value = 32.338836
$ws.Range("O10").Value = 0.1321092878737708
$ws.Range("P10").Value = 0.1321092878737708
$ws.Range("Q10").Value = 35.172174370508
$ws.Range("R10").Value = 316.549569334572
$ws.Range("S10").Value = 0.001845126881051484
$ws.Range("T10").Value = 0.001845126881051483
# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 3.262842333333333
$ws.Range("H11").Value = 9.788527
$ws.Range("I11").Value = 0.01396667040408609
$ws.Range("J11").Value = 0.01396667040408609
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 38.54369466666667
$ws.Range("N11").Value = 115.631084
$ws.Range("O11").Value = 0.4723713668393066
$ws.Range("P11").Value = 0.4723713668393065
$ws.Range("Q11").Value = 125.7619986414742
$ws.Range("R11").Value = 1131.857987773268
$ws.Range("S11").Value = 0.006597455188972236
$ws.Range("T11").Value = 0.006597455188972234
# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 3.262842333333333
$ws.Range("H12").Value = 9.788527
$ws.Range("I12").Value = 0.01396667040408609
$ws.Range("J12").Value = 0.01396667040408609
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 12.62567333333333
$ws.Range("N12").Value = 37.87702
$ws.Range("O12").Value = 0.1547336502458089
$ws.Range("P12").Value = 0.1547336502458089
$ws.Range("Q12").Value = 41.19558143883778
$ws.Range("R12").Value = 370.76023294954
$ws.Range("S12").Value = 0.002161113893404348
$ws.Range("T12").Value = 0.002161113893404347
# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 3.262842333333333
$ws.Range("H13").Value = 9.788527
$ws.Range("I13").Value = 0.01396667040408609
$ws.Range("J13").Value = 0.01396667040408609
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 19.64719066666667
$ws.Range("N13").Value = 58.94157200000001
$ws.Range("O13").Value = 0.2407856950411137
$ws.Range("P13").Value = 0.2407856950411137
$ws.Range("Q13").Value = 64.10568543827156
$ws.Range("R13").Value = 576.9511689444441
$ws.Range("S13").Value = 0.003362974440658022
$ws.Range("T13").Value = 0.00336297444065802
# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 10.83271733333333
$ws.Range("H14").Value = 32.498152
$ws.Range("I14").Value = 0.04636969155071965
$ws.Range("J14").Value = 0.04636969155071963
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 10.779612
$ws.Range("N14").Value = 32.338836
$ws.Range("O14").Value = 0.1321092878737708
$ws.Range("P14").Value = 0.1321092878737708
$ws.Range("Q14").Value = 116.772489759008
$ws.Range("R14").Value = 1050.952407831072
$ws.Range("S14").Value = 0.00612586692969198
$ws.Range("T14").Value = 0.006125866929691978
# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 10.83271733333333
$ws.Range("H15").Value = 32.498152
$ws.Range("I15").Value = 0.04636969155071965
$ws.Range("J15").Value = 0.04636969155071963
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 38.54369466666667
$ws.Range("N15").Value = 115.631084
$ws.Range("O15").Value = 0.4723713668393066
$ws.Range("P15").Value = 0.4723713668393065
$ws.Range("Q15").Value = 417.5329493063076
$ws.Range("R15").Value = 3757.796543756768
$ws.Range("S15").Value = 0.02190371457773049
$ws.Range("T15").Value = 0.02190371457773048
# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 10.83271733333333
$ws.Range("H16").Value = 32.498152
$ws.Range("I16").Value = 0.04636969155071965
$ws.Range("J16").Value = 0.04636969155071963
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 12.62567333333333
$ws.Range("N16").Value = 37.87702
$ws.Range("O16").Value = 0.1547336502458089
$ws.Range("P16").Value = 0.1547336502458089
$ws.Range("Q16").Value = 136.7703503630044
$ws.Range("R16").Value = 1230.93315326704
$ws.Range("S16").Value = 0.007174951634415096
$ws.Range("T16").Value = 0.007174951634415092
# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 10.83271733333333
$ws.Range("H17").Value = 32.498152
$ws.Range("I17").Value = 0.04636969155071965
$ws.Range("J17").Value = 0.04636969155071963
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 19.64719066666667
$ws.Range("N17").Value = 58.94157200000001
$ws.Range("O17").Value = 0.2407856950411137
$ws.Range("P17").Value = 0.2407856950411137
$ws.Range("Q17").Value = 212.8324628861049
$ws.Range("R17").Value = 1915.492165974944
$ws.Range("S17").Value = 0.01116515840888209
$ws.Range("T17").Value = 0.01116515840888208
